$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'52.473.04"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "'  -13.11%  "
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.Value = "'2.326.81"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'  -19.53%  "
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.Value = "'0.998"
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = "'  -0.23%  "
$cell.Style = "Normal"

$cell = $ws.Range("E5")
$cell.Value = "'  -16.58%  "
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.Value = "'121.39"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "'  -14.16%  "
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.Value = "'0.997"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = "'  -0.33%  "
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.Value = "'  -14.90%  "
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.Value = "'2.310.99"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = "'  -20.31%  "
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.Value = "'5.20"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'  -13.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.Value = "'0.0888"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'  -17.23%  "
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.Value = "'0.305"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = "'  -14.76%  "
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.Value = "'  -5.34%  "
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.Value = "'52.534.41"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'  -13.12%  "
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.Value = "'18.94"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'  -16.55%  "
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.Value = "'  -15.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.Value = "'2.339.50"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'  -19.54%  "
$cell.Style = "Normal"

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D18")
$cell.Value = "'3.96"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'  -20.30%  "
$cell.Style = "Normal"

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell = $ws.Range("D19")
$cell.Value = "'300.67"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'  -15.70%  "
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.Value = "'8.97"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'  -22.77%  "
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.Value = "'0.997"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "'  -0.27%  "
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.Value = "'  -1.46%  "
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.Value = "'5.17"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "'  -22.15%  "
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.Value = "'53.61"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "'  -17.11%  "
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.Value = "'0.365"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'  -19.18%  "
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.Value = "'0.146"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "'  -19.31%  "
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.Value = "'6.99"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'  -10.76%  "
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.Value = "'0.997"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.Value = "'  -0.27%  "
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.Value = "'0.0₃0676"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = "'  -19.21%  "
$cell.Style = "Normal"

$cell = $ws.Range("D30")
$cell.Value = "'143.02"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = "'  -4.67%  "
$cell.Style = "Normal"

$cell = $ws.Range("D31")
$cell.Value = "'17.00"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = "'  -13.58%  "
$cell.Style = "Normal"

$cell = $ws.Range("E32")
$cell.Value = "'  -20.39%  "
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.Value = "'4.75"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = "'  -14.55%  "
$cell.Style = "Normal"

$cell = $ws.Range("D34")
$cell.Value = "'3.51"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.Value = "'  -18.76%  "
$cell.Style = "Normal"

$cell = $ws.Range("D35")
$cell.Value = "'0.825"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = "'  -17.09%  "
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = "'  -16.00%  "
$cell.Style = "Normal"

$cell = $ws.Range("D37")
$cell.Value = "'0.991"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = "'  -0.77%  "
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.Value = "'31.81"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = "'  -15.66%  "
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.Value = "'  -1.44%  "
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.Value = "'3.16"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "'  -14.53%  "
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.Value = "'0.0504"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "'  -13.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.Value = "'  -17.31%  "
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.Value = "'1.909.65"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'  -16.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.Value = "'0.521"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'  -19.43%  "
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.Value = "'0.0208"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "'  -12.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.Value = "'0.0827"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'  -9.65%  "
$cell.Style = "Normal"

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D47")
$cell.Value = "'4.04"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "'  -18.22%  "
$cell.Style = "Normal"

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D48")
$cell.Value = "'15.76"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "'  -22.00%  "
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.Value = "'  -5.79%  "
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.Value = "'15.20"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "'  -16.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.Value = "'4.45"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "'  -13.70%  "
$cell.Style = "Normal"
